$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# New caption textbox for the forgotten picture, matching the layout used
# by the other "<filename>.png" caption textboxes in this deck (e.g. the
# "pedestrian-top-forces.png" box on slide 5 which sits at the same
# position/size: left=0, top=11668 EMU, 3048000 x 369332 EMU).
$shp = $s.Shapes.AddTextbox(1, 0, 0.91874015748, 240, 29.0812598425)
$shp.Name = "TextBox 47"
$shp.Fill.Visible = $false
$shp.TextFrame.WordWrap = $true
$shp.TextFrame.AutoSize = 1
$shp.TextFrame.TextRange.Text = "pedestrian-internal-forces.png"
